$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at X (shifts old X/Y/Z -> Y/Z/AA) and merge the
# new column's header cells (rows 1-3) like the other header columns.
$ws.Range("X1:X3").EntireColumn.Insert()
$ws.Range("X1:X3").Merge()

# New header text for the inserted "Institution" column.
$ws.Range("X1").Value = "Institution"

# New trailing "(ab)" label cell for the row-4 letter sequence
# (was Z4 = "(aa)", now AA4 is the next free cell after the insert).
$ws.Range("AA4").Value = "(ab)"

# Restore/update the view state (selection + scroll position) to match.
$ws.Range("AA4").Select()
$excel.ActiveWindow.ScrollColumn = 2
